{"js": "// Adds two new bullet points after the paragraph that ends the\n// \"Bei der Briefanschrift werden immer alle verf\u00fcgbaren Anreden aller\n// bekannten Titel verwendet\" list, keeping the same \"Listenabsatz\" style\n// and the same numbered-list membership as their neighbours.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"Bei der Briefanschrift werden immer alle verf\u00fcgbaren Anreden aller bekannten Titel verwendet\";\n\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === anchorText) {\n    anchor = p;\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph for the new bullet points.\");\n}\n\n// Capture the list this paragraph belongs to, so the new paragraphs can be\n// attached to the very same numbering list (numId) at the same level.\nconst list = anchor.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\nconst firstNew = anchor.insertParagraph(\n  \"Die Briefanschrift wird, wenn erkannt in Deutsch verfasst. Ansonsten wird das englische Pr\u00e4fix \u201eDear\u201c verwendet.\",\n  Word.InsertLocation.after\n);\nfirstNew.style = \"Listenabsatz\";\nfirstNew.attachToList(listId, 0);\n\nconst secondNew = firstNew.insertParagraph(\n  \"Wenn die Briefanschrift nicht bestimmbar ist, wird der englische Standard \u201eDear Sirs\u201c verwendet.\",\n  Word.InsertLocation.after\n);\nsecondNew.style = \"Listenabsatz\";\nsecondNew.attachToList(listId, 0);\n\nawait context.sync();\n", "ps1": "# Adds two new bullet points after the paragraph that ends the\n# \"Bei der Briefanschrift werden immer alle verf\u00fcgbaren Anreden aller\n# bekannten Titel verwendet\" list, keeping the same \"Listenabsatz\" style\n# and list numbering as the surrounding bullets (InsertParagraphAfter on a\n# list paragraph's Range inherits the paragraph's style/numbering, exactly\n# like pressing Enter at the end of that list item in Word).\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Bei der Briefanschrift werden immer alle verf\u00fcgbaren Anreden aller bekannten Titel verwendet\"\n\n$target = $null\n$targetIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n    $i = $i + 1\n}\nif ($target -eq $null) {\n    throw \"Could not find the anchor paragraph for the new bullet points.\"\n}\n\n# First new bullet, inserted right after the anchor paragraph.\n$target.Range.InsertParagraphAfter()\n$firstNew = $d.Paragraphs.Item($targetIndex + 1)\n$firstNew.Range.Text = \"Die Briefanschrift wird, wenn erkannt in Deutsch verfasst. Ansonsten wird das englische Pr\u00e4fix \u201eDear\u201c verwendet.\"\n\n# Second new bullet, inserted right after the first new one.\n$firstNew.Range.InsertParagraphAfter()\n$secondNew = $d.Paragraphs.Item($targetIndex + 2)\n$secondNew.Range.Text = \"Wenn die Briefanschrift nicht bestimmbar ist, wird der englische Standard \u201eDear Sirs\u201c verwendet.\"\n"}
